# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the two trailing worker rows that no longer belong in this
#        report (ALEGUIS ENRIQUE DE AVILA MARTELO / DIEGO ARMANDO GONZALEZ
#        DEVIA). They sit right after the GINA MARGARITA BARRIOS FRANCESCH
#        block, at rows 62 and 63. Delete from the bottom up so row numbers
#        above stay valid while deleting.
$ws.Rows(63).Delete()
$ws.Rows(62).Delete()

# --- 2. The remaining 46 detail rows (16-61) for GINA MARGARITA BARRIOS
#        FRANCESCH were listed newest-period-first (2208 down to 1811).
#        The new workbook lists them oldest-period-first (1811 up to 2208).
#        Capture the "Periodo Mora" (col E) and "Valor Mora" (col F) values,
#        reverse their order, and write them back in place - every other
#        column (B, C, D, G, H, I, J) is identical for all of these rows so
#        it does not need to move.
$firstRow = 16
$lastRow = 61

$periodos = @()
$valores = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodos += , $ws.Cells.Item($r, 5).Value()
    $valores += , $ws.Cells.Item($r, 6).Value()
}

$periodos = $periodos[-1..-($periodos.Length)]
$valores = $valores[-1..-($valores.Length)]

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 5).Value = $periodos[$i]
    $ws.Cells.Item($r, 6).Value = $valores[$i]
}

# --- 3. Refresh the summary figures at the top of the report.
$ws.Range("E11").Value = 1475357
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 46
